$d = $word.ActiveDocument

$replacements = @(
    @{old = "67×65="; new = "81×18="},
    @{old = "25×24="; new = "54×69="},
    @{old = "93×12="; new = "30×90="},
    @{old = "42×89="; new = "56×92="},
    @{old = "91×22="; new = "99×74="},
    @{old = "73×24="; new = "82×97="},
    @{old = "95×29="; new = "62×65="},
    @{old = "35×27="; new = "68×18="},
    @{old = "39×19="; new = "98×25="},
    @{old = "90×62="; new = "94×80="},
    @{old = "16×36="; new = "92×66="},
    @{old = "33×53="; new = "23×62="},
    @{old = "59×82="; new = "74×35="},
    @{old = "80×16="; new = "21×12="},
    @{old = "38×43="; new = "60×36="},
    @{old = "29×94="; new = "42×68="},
    @{old = "86×71="; new = "45×29="},
    @{old = "13×39="; new = "78×53="},
    @{old = "93×86="; new = "65×67="},
    @{old = "26×43="; new = "22×14="},
    @{old = "64×98="; new = "92×73="},
    @{old = "43×32="; new = "70×94="},
    @{old = "64×25="; new = "87×42="},
    @{old = "77×77="; new = "68×51="},
    @{old = "18×42="; new = "76×38="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
